# Applies the "gh-pages output regenerated" update to 苏州-漫展信息.xlsx
# Target sheets: "展览" (index 1) and "全部类型" (index 4) -- both hold the
# same exhibition table and receive the identical update.

$wb = $excel.ActiveWorkbook

function Set-EventRow {
    param($ws, $r, $idx, $date, $name, $place, $timeRange, $want, $price, $hasStage, $cover)

    $ws.Cells.Item($r, 1).Value2 = $idx

    # Dates are stored as literal text like "2024-01-21" -- a leading
    # apostrophe keeps the engine from coercing the date-shaped text to a
    # serial date number (mirrors typing '2024-01-21 into Excel).
    $ws.Cells.Item($r, 2).Value2 = "'" + $date

    if ($name -ne $null -and $name -ne "") {
        $ws.Cells.Item($r, 3).Value2 = $name
    } else {
        $ws.Cells.Item($r, 3).ClearContents()
    }

    $ws.Cells.Item($r, 4).Value2 = $place
    $ws.Cells.Item($r, 5).Value2 = $timeRange
    $ws.Cells.Item($r, 6).Value2 = $want

    # Min-price column is text ("60", "49", ...) except for the two rows
    # that hold a status word instead ("已售罄"/"预售中"); only force the
    # quote-prefix when the value actually looks numeric.
    if ($price -match "^[0-9]+$") {
        $ws.Cells.Item($r, 7).Value2 = "'" + $price
    } else {
        $ws.Cells.Item($r, 7).Value2 = $price
    }

    $ws.Cells.Item($r, 8).Value2 = $hasStage
    $ws.Cells.Item($r, 9).ClearContents()
    $ws.Cells.Item($r, 10).Value2 = $cover
}

$sheetIndexes = @(1, 4)
foreach ($si in $sheetIndexes) {
    $ws = $wb.Worksheets.Item($si)

    # Row 17 is brand new -- clone formatting from row 16 (bold/border/
    # centered style used by every row-index cell in column A) before
    # writing into it so the new row matches the existing table style.
    $ws.Cells.Item(16, 1).Copy()
    $ws.Cells.Item(17, 1).PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    Set-EventRow $ws 2 1 "2024-01-21" "苏州·Good Jump ACG迎新特别篇X动漫品牌博览会" "金山南路288号 广电国际会展中心" "2024.01.21 10:30-01.21 17:00" 1900 "60" $true "//i1.hdslb.com/bfs/openplatform/202312/vtGcfnyc1703060683812.jpeg"
    Set-EventRow $ws 3 2 "2024-01-27" "苏州·第五届次元鹿角动漫游戏展" "绿地大道258号游站未来城2号楼 魔之塔" "2024.01.27 10:00-01.27 19:00" 269 "49" $false "//i2.hdslb.com/bfs/openplatform/202311/Z7mV6VXN1701160508967.jpeg"
    Set-EventRow $ws 4 3 "2024-01-28" "昆山·“不是！你有病吧！”主题展" "白塔东路60号(近平江路) 苏州书香府邸平江府" "2024.01.28 10:00-01.28 21:00" 255 "58" $false "//i0.hdslb.com/bfs/openplatform/202311/5AgvDWGQ1700817845950.jpeg"
    Set-EventRow $ws 5 4 "2024-02-03" "苏州·世纪幻想动漫游戏展" "苏州大道东688号 苏州国际博览中心" "2024.02.03 09:30-02.04 17:00" 8586 "60" $false "//i0.hdslb.com/bfs/openplatform/202401/aDe3s9MS1705479547745.jpeg"
    Set-EventRow $ws 6 5 "2024-02-03" "苏州·第十七届 I COME ACG  动漫品牌博览会" "金山南路288号木渎影视城F2 苏州广电国际会展中心" "2024.02.03 10:00-02.03 20:00" 10129 "25" $false "//i2.hdslb.com/bfs/openplatform/202401/IkyhIHPT1704352086775.jpeg"
    Set-EventRow $ws 7 6 "2024-02-04" $null "苏州大道东688号 苏州国际博览中心" "2024.02.04 09:30-02.04 17:00" 577 "已售罄" $false "//i1.hdslb.com/bfs/openplatform/202401/bHsHJ3f21704186294427.jpeg"
    Set-EventRow $ws 8 7 "2024-02-14" "苏州·国风宠物-cosplay展" "润元路润南巷172号,地铁二号线陆慕站东200米,近市旅游换乘中心北100米 斐利酒店" "2024.02.14 10:00-02.14 16:00" 21 "49" $true "//i2.hdslb.com/bfs/openplatform/202401/oWbVnOjD1704445446390.jpeg"
    Set-EventRow $ws 9 8 "2024-02-14" "苏州·梦幻岛 国乙主题文化展（日夜场） 梦幻岛之约3.0" "常熟国际展览中心 国际展览中心" "2024.02.14 09:00-02.15 17:30" 625 "55" $false "//i1.hdslb.com/bfs/openplatform/202401/VHHzVjad1704438989848.jpeg"
    Set-EventRow $ws 10 9 "2024-02-14" "苏州·第一届寒假动漫展宅舞比赛-CF01" "虞山北路258号 星程酒店(长江路店)" "2024.02.14 09:00-02.14 21:00" 86 "50" $false "//i2.hdslb.com/bfs/openplatform/202312/oPrKUOby1703664065719.jpeg"
    Set-EventRow $ws 11 10 "2024-02-16" "常熟·CDW·动漫展02" "金山南路288号 广电国际会展中心" "2024.02.16 10:00-02.17 17:00" 9291 "60" $false "//i2.hdslb.com/bfs/openplatform/202312/C3P0Encm1701659824998.jpeg"
    Set-EventRow $ws 12 11 "2024-02-25" "苏州·绘时国乙1.0-秩序之外" "清禾路886号 尹山湖大剧院" "2024.02.25 10:00-02.25 17:00" 2401 "68" $true "//i1.hdslb.com/bfs/openplatform/202401/tqrMA6qB1704787264871.jpeg"
    Set-EventRow $ws 13 12 "2024-03-08" "苏州·TCD国潮动漫游戏嘉年华吴磊内场" "木渎金山南路288号 苏州国际影视娱乐城" "2024.03.08 09:00-03.10 17:30" 15 "65" $true "//i2.hdslb.com/bfs/openplatform/202401/Rfd9PcBN1704781416369.jpeg"
    Set-EventRow $ws 14 13 "2024-04-13" "常熟·漫魂动漫游戏展01" "石路步行街永福桥浜15号 银河广场" "2024.04.13 13:30-04.13 20:00" 45 "78" $false "//i0.hdslb.com/bfs/openplatform/202401/SjKfDxBh1705041298410.jpeg"
    Set-EventRow $ws 15 14 "2024-04-21" "苏州.第二届THO 赤维极陵" "清禾路888号2号楼3楼 格莱美婚礼宴会中心" "2024.04.21 10:00-04.21 21:00" 337 "48" $true "//i0.hdslb.com/bfs/openplatform/202312/X0PZ3YhH1703822037665.jpeg"
    Set-EventRow $ws 16 15 "2024-05-01" "昆山·第十二届理想乡动漫游戏展" "花桥经济开发区绿地大道1598号 花桥国际博览中心" "2024.05.01 10:00-05.03 17:00" 10606 "65" $true "//i2.hdslb.com/bfs/openplatform/202312/lau3mW031702535438289.jpeg"
    Set-EventRow $ws 17 16 "2024-05-01" "苏州·TCD国潮动漫游戏嘉年华" "花桥经济开发区绿地大道1598号 花桥国际博览中心" "2024.05.01 10:00-05.03 17:00" 10295 "预售中" $true "//i2.hdslb.com/bfs/openplatform/202310/9xMTQMlg1696736126094.png"
}
